$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '23.460.92'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +1.13%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.638.60'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +2.23%  '
$ws.Range('E4').Value = '  +0.06%  '
$ws.Range('E5').Value = '  +0.04%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '306.30'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.93%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.3753'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -0.76%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '51.99'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -0.11%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.3637'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +0.45%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '1.260'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -0.85%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.08139'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +0.24%  '
$ws.Range('E12').Value = '  +0.04%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '22.94'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +0.59%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.623'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +0.42%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.00001275'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +2.41%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '7.365'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -0.80%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '1.637.98'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +2.21%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '94.50'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +0.41%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.06908'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +0.40%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '18.17'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +0.41%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.534'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -0.27%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.9999'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.03%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '23.472.59'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +1.16%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '12.76'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -1.76%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '3.084'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +3.49%  '
$ws.Range('E26').Value = '  +0.88%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '21.24'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -0.12%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '151.15'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +1.19%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '5.324'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +1.32%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '137.13'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +2.35%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '2.307'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -2.82%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.817.91'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +2.21%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '6.757'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -0.45%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.9602'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -1.10%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.02831'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +3.97%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '10.34'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +0.30%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.07293'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -3.03%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.2527'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +0.93%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.08823'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +0.19%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '6.117'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +0.65%  '
$ws.Range('E41').Value = '  +1.27%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.7089'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -0.40%  '
$ws.Range('B43').Value = 'Aptos'
$ws.Range('C43').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '12.46'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -0.45%  '
$ws.Range('B44').Value = 'EnergySwap'
$ws.Range('C44').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '16.20'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +3.49%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.6549'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +0.16%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.336'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +0.94%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.9999'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +0.32%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '4.012'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -0.11%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.07969'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +0.11%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '128.63'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -2.83%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.205'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +0.30%  '
